$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.17622703276409
$ws.Range("C2").Value = 8.448252368653622
$ws.Range("D2").Value = 3.852438476026343
$ws.Range("E2").Value = 11.87252651299297
$ws.Range("F2").Value = 20.85295102378897
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 17.4299331331869
$ws.Range("M2").Value = 14.15154519726645
$ws.Range("N2").Value = 16.33347107265294
$ws.Range("O2").Value = 18.4545401647835
$ws.Range("B3").Value = 10.61141056364527
$ws.Range("C3").Value = 8.010231711293601
$ws.Range("D3").Value = 3.805593125082941
$ws.Range("E3").Value = 11.75994438436761
$ws.Range("F3").Value = 20.77413510526166
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 17.51116580657236
$ws.Range("M3").Value = 13.87218533177993
$ws.Range("N3").Value = 16.38461062119998
$ws.Range("O3").Value = 18.46272471924573
$ws.Range("B4").Value = 10.24956465618201
$ws.Range("C4").Value = 7.727085927984437
$ws.Range("D4").Value = 3.776170009103248
$ws.Range("E4").Value = 11.69489938987268
$ws.Range("F4").Value = 20.73284152840582
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 17.56551027831303
$ws.Range("M4").Value = 13.70101718267319
$ws.Range("N4").Value = 16.41776382012985
$ws.Range("O4").Value = 18.47334203089336
$ws.Range("B5").Value = 10.09849470632515
$ws.Range("C5").Value = 7.608197006260048
$ws.Range("D5").Value = 3.76402063996745
$ws.Range("E5").Value = 11.66944737282193
$ws.Range("F5").Value = 20.71781252560335
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 17.5887757252453
$ws.Range("M5").Value = 13.63145501701302
$ws.Range("N5").Value = 16.43171584331053
$ws.Range("O5").Value = 18.4790719967259
$ws.Range("B6").Value = 10.07319683661221
$ws.Range("C6").Value = 7.588246105740095
$ws.Range("D6").Value = 3.761993829930942
$ws.Range("E6").Value = 11.66528553646779
$ws.Range("F6").Value = 20.71542592454844
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 17.59270644972904
$ws.Range("M6").Value = 13.61991870714803
$ws.Range("N6").Value = 16.43405927928135
$ws.Range("O6").Value = 18.48010810958367
$ws.Range("B7").Value = 10.24754166594819
$ws.Range("C7").Value = 7.725496644043808
$ws.Range("D7").Value = 3.776006794119042
$ws.Range("E7").Value = 11.69455183079312
$ws.Range("F7").Value = 20.73263154478496
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 17.56581951583456
$ws.Range("M7").Value = 13.70007813978186
$ws.Range("N7").Value = 16.41795019180918
$ws.Range("O7").Value = 18.47341362966682
$ws.Range("B8").Value = 10.9846937104618
$ws.Range("C8").Value = 8.300205396035716
$ws.Range("D8").Value = 3.836427139027142
$ws.Range("E8").Value = 11.83287889940988
$ws.Range("F8").Value = 20.82431050332252
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 17.45701275459138
$ws.Range("M8").Value = 14.05520650660509
$ws.Range("N8").Value = 16.35074065412253
$ws.Range("O8").Value = 18.45620017442541
$ws.Range("B9").Value = 12.3048017799985
$ws.Range("C9").Value = 9.312272240193963
$ws.Range("D9").Value = 3.949389528930688
$ws.Range("E9").Value = 12.13515722279186
$ws.Range("F9").Value = 21.05971489511903
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 17.27925375878512
$ws.Range("M9").Value = 14.75018965945458
$ws.Range("N9").Value = 16.23281207741392
$ws.Range("O9").Value = 18.46689838648336
$ws.Range("B10").Value = 13.19169841628281
$ws.Range("C10").Value = 9.983448307774093
$ws.Range("D10").Value = 4.028648332276102
$ws.Range("E10").Value = 12.37417730824797
$ws.Range("F10").Value = 21.26543920845316
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 17.17058806886872
$ws.Range("M10").Value = 15.25440062742455
$ws.Range("N10").Value = 16.15456584690864
$ws.Range("O10").Value = 18.50191149302952
$ws.Range("B11").Value = 13.57615563574218
$ws.Range("C11").Value = 10.27273903397841
$ws.Range("D11").Value = 4.063820199713505
$ws.Range("E11").Value = 12.48613609818137
$ws.Range("F11").Value = 21.36584941934503
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 17.12596383555277
$ws.Range("M11").Value = 15.48129520752808
$ws.Range("N11").Value = 16.12078046575501
$ws.Range("O11").Value = 18.52372566244056
$ws.Range("B12").Value = 13.71894466047142
$ws.Range("C12").Value = 10.37996137389608
$ws.Range("D12").Value = 4.077005891149779
$ws.Range("E12").Value = 12.52895625092649
$ws.Range("F12").Value = 21.40482522645551
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 17.10976119223364
$ws.Range("M12").Value = 15.56677328528751
$ws.Range("N12").Value = 16.10824611085542
$ws.Range("O12").Value = 18.53282991719347
$ws.Range("B13").Value = 13.6883177705603
$ws.Range("C13").Value = 10.35697283830766
$ws.Range("D13").Value = 4.074172132284549
$ws.Range("E13").Value = 12.51971596567685
$ws.Range("F13").Value = 21.39638920126546
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 17.11321971927098
$ws.Range("M13").Value = 15.54838502997225
$ws.Range("N13").Value = 16.11093408131483
$ws.Range("O13").Value = 18.53083168487971
$ws.Range("B14").Value = 13.58795935021578
$ws.Range("C14").Value = 10.28160699237813
$ws.Range("D14").Value = 4.064907696752655
$ws.Range("E14").Value = 12.48965067785505
$ws.Range("F14").Value = 21.36903706267555
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 17.12461687010614
$ws.Range("M14").Value = 15.48833687491132
$ws.Range("N14").Value = 16.11974406144148
$ws.Range("O14").Value = 18.52445778566263
$ws.Range("B15").Value = 13.52612100863873
$ws.Range("C15").Value = 10.23513982666652
$ws.Range("D15").Value = 4.059215447736457
$ws.Range("E15").Value = 12.4712888026315
$ws.Range("F15").Value = 21.35240625849981
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 17.1316886624849
$ws.Range("M15").Value = 15.4514956113307
$ws.Range("N15").Value = 16.12517419363795
$ws.Range("O15").Value = 18.52066336239815
$ws.Range("B16").Value = 13.16618513648154
$ws.Range("C16").Value = 9.964218021207701
$ws.Range("D16").Value = 4.026331465215085
$ws.Range("E16").Value = 12.36692205285463
$ws.Range("F16").Value = 21.2590120910997
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 17.17360147919116
$ws.Range("M16").Value = 15.23951527139815
$ws.Range("N16").Value = 16.15681014439126
$ws.Range("O16").Value = 18.50060412554834
$ws.Range("B17").Value = 12.94046039803226
$ws.Range("C17").Value = 9.793894094314936
$ws.Range("D17").Value = 4.005927538100794
$ws.Range("E17").Value = 12.30369461362844
$ws.Range("F17").Value = 21.20344607562665
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 17.2005481295515
$ws.Range("M17").Value = 15.10877572963395
$ws.Range("N17").Value = 16.17668063482697
$ws.Range("O17").Value = 18.48980448738361
$ws.Range("B18").Value = 12.80884445391968
$ws.Range("C18").Value = 9.694420261452674
$ws.Range("D18").Value = 3.994108943174727
$ws.Range("E18").Value = 12.26763410885659
$ws.Range("F18").Value = 21.1721302596117
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 17.2164994667725
$ws.Range("M18").Value = 15.03334984215863
$ws.Range("N18").Value = 16.18827994983419
$ws.Range("O18").Value = 18.48414712445035
$ws.Range("B19").Value = 12.76397716079547
$ws.Range("C19").Value = 9.660481841790984
$ws.Range("D19").Value = 3.990093339696112
$ws.Range("E19").Value = 12.25547845482135
$ws.Range("F19").Value = 21.16163876186655
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 17.22197787768594
$ws.Range("M19").Value = 15.00777548948636
$ws.Range("N19").Value = 16.19223655629738
$ws.Range("O19").Value = 18.48232690742185
$ws.Range("B20").Value = 12.96467443719921
$ws.Range("C20").Value = 9.812181606562135
$ws.Range("D20").Value = 4.00810818499475
$ws.Range("E20").Value = 12.31039387674396
$ws.Range("F20").Value = 21.20929469352026
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 17.19763276726139
$ws.Range("M20").Value = 15.12271741080289
$ws.Range("N20").Value = 16.1745477635585
$ws.Range("O20").Value = 18.49089677560295
$ws.Range("B21").Value = 13.61751340535383
$ws.Range("C21").Value = 10.30380702810201
$ws.Range("D21").Value = 4.067632546655165
$ws.Range("E21").Value = 12.49847040478008
$ws.Range("F21").Value = 21.37704543452549
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 17.12125033972368
$ws.Range("M21").Value = 15.50598713183449
$ws.Range("N21").Value = 16.11714932139387
$ws.Range("O21").Value = 18.52630708402744
$ws.Range("B22").Value = 14.02785902395909
$ws.Range("C22").Value = 10.61155010428534
$ws.Range("D22").Value = 4.105756505171667
$ws.Range("E22").Value = 12.62384034453096
$ws.Range("F22").Value = 21.49221814268589
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 17.07538631251707
$ws.Range("M22").Value = 15.75386339120506
$ws.Range("N22").Value = 16.08114799037181
$ws.Range("O22").Value = 18.55436545671366
$ws.Range("B23").Value = 13.810361676519
$ws.Range("C23").Value = 10.44854833299804
$ws.Range("D23").Value = 4.085482232387184
$ws.Range("E23").Value = 12.55671725300024
$ws.Range("F23").Value = 21.43025163653849
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 17.09949229214143
$ws.Range("M23").Value = 15.62183344261981
$ws.Range("N23").Value = 16.10022447652584
$ws.Range("O23").Value = 18.53894158636497
$ws.Range("B24").Value = 12.9537330080197
$ws.Range("C24").Value = 9.803918655917467
$ws.Range("D24").Value = 4.007122588489841
$ws.Range("E24").Value = 12.30736423416208
$ws.Range("F24").Value = 21.20664856931424
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 17.19894937198984
$ws.Range("M24").Value = 15.11641518613716
$ws.Range("N24").Value = 16.175511487894
$ws.Range("O24").Value = 18.49040123400422
$ws.Range("B25").Value = 11.96185823515197
$ws.Range("C25").Value = 9.05105970946682
$ws.Range("D25").Value = 3.91946015096175
$ws.Range("E25").Value = 12.05025863478769
$ws.Range("F25").Value = 20.9901840011886
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 17.32350665063144
$ws.Range("M25").Value = 14.56291480160231
$ws.Range("N25").Value = 16.23281207741392
$ws.Range("O25").Value = 18.45923558686932
